$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new "status" column (I) for the data rows: row 2 passed,
# rows 3-7 failed.
$ws.Range("I2").Value = "PASSED"
$ws.Range("I3").Value = "FAILED"
$ws.Range("I4").Value = "FAILED"
$ws.Range("I5").Value = "FAILED"
$ws.Range("I6").Value = "FAILED"
$ws.Range("I7").Value = "FAILED"
